# Update the roster table on Sheet1: rows 4-6 and 8-16 (1-indexed Excel rows)
# with new Player / Position / Team values, reflecting the re-uploaded roster.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A4").Value = "Collin Sexton"
$ws.Range("B4").Value = "PG,SG"
$ws.Range("C4").Value = "Utah Jazz"

$ws.Range("A5").Value = "Malik Monk"
$ws.Range("B5").Value = "PG,SG,SF"
$ws.Range("C5").Value = "Sacramento Kings"

$ws.Range("A6").Value = "Devin Vassell"
$ws.Range("B6").Value = "SG,SF"
$ws.Range("C6").Value = "San Antonio Spurs"

$ws.Range("A8").Value = "Tim Hardaway Jr."
$ws.Range("B8").Value = "SG,SF"
$ws.Range("C8").Value = "Detroit Pistons"

$ws.Range("A9").Value = "Onyeka Okongwu"
$ws.Range("B9").Value = "PF,C"
$ws.Range("C9").Value = "Atlanta Hawks"

$ws.Range("A12").Value = "Naz Reid"
$ws.Range("B12").Value = "PF,C"
$ws.Range("C12").Value = "Minnesota Timberwolves"

$ws.Range("A13").Value = "Precious Achiuwa"
$ws.Range("B13").Value = "PF,C"
$ws.Range("C13").Value = "New York Knicks"

$ws.Range("A14").Value = "Derrick White"
$ws.Range("B14").Value = "PG,SG"
$ws.Range("C14").Value = "Boston Celtics"

$ws.Range("A15").Value = "Coby White"
$ws.Range("B15").Value = "PG,SG"
$ws.Range("C15").Value = "Chicago Bulls"

$ws.Range("A16").Value = "Duncan Robinson"
$ws.Range("B16").Value = "SG,SF"
$ws.Range("C16").Value = "Miami Heat"
